# Update consolidated log file
# The "Latency" block (header + Threads/Amazon/Google table) that used to
# live at rows 163, 165-172 on the "Test2" sheet is relocated up to rows
# 112, 114-121, and the "Google" (column C) latency figures are refreshed
# with new measurements.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test2")

# --- Move the block from its old location (rows 163-172) to the new one
#     (rows 112-121), preserving values + formatting via Cut/Paste. ---
$ws.Range("A163").Cut($ws.Range("A112"))

$ws.Range("A165:C165").Cut($ws.Range("A114:C114"))
$ws.Range("A166:C166").Cut($ws.Range("A115:C115"))
$ws.Range("A167:C167").Cut($ws.Range("A116:C116"))
$ws.Range("A168:C168").Cut($ws.Range("A117:C117"))
$ws.Range("A169:C169").Cut($ws.Range("A118:C118"))
$ws.Range("A170:C170").Cut($ws.Range("A119:C119"))
$ws.Range("A171:C171").Cut($ws.Range("A120:C120"))
$ws.Range("A172:C172").Cut($ws.Range("A121:C121"))

# --- Remove any leftover (now-empty) formatting stubs at the old rows so
#     they no longer appear in the sheet at all. ---
$ws.Range("A163:C172").Clear()

# --- Refresh the "Google" (column C) latency values with the new data ---
$ws.Range("C115").Value = 67.79
$ws.Range("C116").Value = 108.4
$ws.Range("C117").Value = 214.54
$ws.Range("C118").Value = 453.34
$ws.Range("C119").Value = 1009.34
$ws.Range("C120").Value = 2752.25
$ws.Range("C121").Value = 8556.25
